$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading_percent values for columns B,C,E,F,G,H,I,K,M across rows 2-25
# (columns D,J,L,N,O are unchanged zeros)
$data = @(
    @(2, 8.461437329925868, 4.460189247900223, 12.485333666567, 16.86991607391233, 34.37797511545209, 15.67023792101732, 22.8506990301772, 8.66227387128604, 13.74349141120732),
    @(3, 8.195696921778865, 4.220086105404775, 12.2675410576749, 15.89584955866808, 34.4177185865039, 15.72230577637472, 22.94046083773067, 8.492351832921162, 13.58362259603142),
    @(4, 8.029764301248322, 4.064587014020219, 12.13613561859007, 15.26997757108491, 34.45406274580446, 15.75717084808058, 23.00044431928093, 8.388097393767175, 13.48806807113769),
    @(5, 7.961561047474669, 3.99921506351468, 12.08324729590299, 15.00819731993403, 34.47186061370317, 15.77210524844285, 23.02610900911351, 8.345693092563481, 13.44982839242391),
    @(6, 7.95020406829768, 3.988239882841035, 12.07450735325043, 14.96433081551593, 34.47499587260494, 15.77462893737493, 23.03044425000672, 8.338658525939424, 13.44352224282314),
    @(7, 8.028846701255398, 4.063713464184683, 12.13541957287557, 15.26647399323137, 34.45429070049706, 15.75736931856617, 23.00078550343651, 8.38752510893303, 13.48754946741245),
    @(8, 8.370448629750166, 4.379091243107006, 12.40980369576312, 16.53996406344768, 34.38919150707174, 15.68758916983225, 22.88063615950065, 8.603704421787633, 13.68785490338737),
    @(9, 9.013899790752559, 4.932621898127261, 12.96287081430115, 19.00274580682531, 34.35688033668593, 15.57378574276794, 22.68381628222929, 9.025663861196781, 14.09931582902228),
    @(10, 9.465138240006826, 5.298889426690892, 13.37381896575403, 20.67494806633232, 34.39194445572837, 15.50430309627489, 22.56307501464825, 9.331141294139536, 14.41023900338817),
    @(11, 9.664777901000724, 5.456633362628354, 13.56083677858474, 21.3917225636224, 34.42076149973965, 15.4757797548635, 22.51337375903658, 9.468466036154984, 14.55296557854448),
    @(12, 9.739496603488753, 5.515085284723461, 13.63159164105272, 21.65686569030329, 34.43352790375457, 15.46542378277855, 22.49530826560354, 9.520180456517734, 14.60714795520214),
    @(13, 9.723444800325543, 5.50255372633324, 13.61635736522162, 21.60004134736742, 34.43069593734599, 15.46763430064972, 22.49916533737574, 9.509056410116491, 14.59547359976387),
    @(14, 9.670943078407067, 5.461467965239588, 13.56665954304012, 21.4136618050453, 34.42177462897016, 15.47491883112537, 22.51187233993712, 9.472726660475637, 14.55742077547963),
    @(15, 9.638667620653917, 5.43613460661339, 13.53620757115731, 21.29868154950795, 34.41655158972982, 15.47943884227365, 22.51975422598922, 9.45043471380966, 14.53412841879645),
    @(16, 9.451971782183417, 5.288401349297818, 13.36159190268347, 20.62722412089977, 34.39032048276083, 15.5062293746147, 22.56642858317032, 9.322129326479418, 14.40093307624133),
    @(17, 9.335944875990327, 5.195493970313222, 13.25443548107108, 20.20408069597325, 34.37752743485582, 15.52345571915605, 22.59640278240834, 9.242960620698128, 14.31951568042829),
    @(18, 9.26868297862975, 5.14122116203298, 13.1928133921654, 19.95656407809801, 34.37138016507168, 15.53365410380771, 22.61413467003848, 9.197273570336021, 14.27281056442718),
    @(19, 9.245821148891361, 5.122702149161791, 13.17195352869681, 19.87204792380568, 34.36950662887026, 15.537156901914, 22.62022267925567, 9.181780362938568, 14.25701984782382),
    @(20, 9.348351144410467, 5.205470585512425, 13.26584183682083, 20.24955283636154, 34.37876391151894, 15.52159189402188, 22.59316107799005, 9.251404309330589, 14.32817023012659),
    @(21, 9.686388529495369, 5.473570690234318, 13.58125934205389, 21.46857628470577, 34.4243446975922, 15.47276709252462, 22.50811945498065, 9.483405789030783, 14.56859454610453),
    @(22, 9.902155728651001, 5.64131394129716, 13.78699280314764, 22.22866616901552, 34.46494102078002, 15.44345295830153, 22.45694410877494, 9.633334620589808, 14.72648946742518),
    @(23, 9.787490256942824, 5.552471807959152, 13.67725067914977, 21.82633154458858, 34.44228455815569, 15.45886042496282, 22.48385304494498, 9.553486454958614, 14.64216443872894),
    @(24, 9.34274399831728, 5.200962833305675, 13.26068506881906, 20.22900810905287, 34.37820113896397, 15.5224336113713, 22.59462509743197, 9.247587449158772, 14.32425718408171),
    @(25, 8.843248616737316, 4.789928140207731, 12.8121181103608, 18.34778573295695, 34.35533492712249, 15.60209758573595, 22.73288619305207, 8.912076435268016, 13.98628914017295)
)

$cols = @("B","C","E","F","G","H","I","K","M")

foreach ($entry in $data) {
    $row = $entry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $value = $entry[$i + 1]
        $ws.Range("$col$row").Value = $value
    }
}

Write-Host "Updated $($data.Length) rows"
